$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 4-13 held the login test case data (TC-004 .. TC-014) which was
# removed; clear their contents (formatting/styles are left untouched).
$ws.Range("A4:BU13").ClearContents()

# Drop every hyperlink except the one on E2, then recreate the E2
# hyperlink (engine only supports bulk-clearing the hyperlink collection),
# preserving its display text, underlying value and style.
$ws.Range("A1").Hyperlinks.Delete()
$h = $ws.Hyperlinks.Add($ws.Range("E2"), "mailto:jamessmith@mailinator.com")
$h.TextToDisplay = "jamessmith@mailinator.com"
$ws.Range("E2").Value2 = "Nitesh"
$ws.Range("E2").Style = "Hyperlink"

# Update the visible view: scroll so column C is left-most and select F2.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("F2").Select()
